$wb = $excel.ActiveWorkbook
$wsTreino = $wb.Worksheets.Item(1)
$wsTeste  = $wb.Worksheets.Item(2)

# --- Column B ("Relevante") classification on the "Teste" sheet ---------
$wsTeste.Cells.Item(1, 2).Value2 = "Relevante"

$classificacoes = @(1,0,0,1,1,1,1,1,1,0,1,0,1,0,0,1,0,0,1,1,1,1,1,1,1,0,0,1,1,0,1,0,0,1,1,0,0,0,1,1,1,0,1,0,1,1,0,1,1,1,1,0,0,1,1,0,1,1,1,0,1,0,1,1,0,0,1,1,0,0,0,0,1,0,1,0,0,1,0,0,1,1,1,1,0,1,0,0,0,0,1,0,0,1,1,0,0,1,1,1,0,1,1,1,1,1,1,1,0,0,1,0,0,1,0,1,0,1,0,1,0,1,0,1,0,1,1,1,1,1,1,0,0,0,1,1,1,0,1,1,1,1,1,0,0,1,1,0,0,0,1,1,1,0,1,0,1,1,0,0,1,0,1,0,1,0,0,0,1,1,0,1,0,1,1,0,1,1,1,1,1,1,1,1,1,0,1,0,1,1,1,1,1,1,0,1,0,0,1,1)

for ($i = 0; $i -lt $classificacoes.Length; $i++) {
    $row = $i + 2
    $wsTeste.Cells.Item($row, 2).Value2 = $classificacoes[$i]
}

# one extra trailing row with only a value in column B (no text in A)
$wsTeste.Cells.Item(202, 2).Value2 = 1

# --- cosmetic bits that mirror the "Treinamento" sheet ------------------
$wsTeste.Columns.Item(1).ColumnWidth = 254.67
$wsTeste.Range("A20").Font.Underline = 2

$wsTeste.PageSetup.PaperSize = 9
$wsTeste.PageSetup.Orientation = 1

# --- selection / active-sheet bookkeeping -------------------------------
$wsTreino.Activate()
$wsTreino.Range("A152").Select()

$wsTeste.Activate()
$wsTeste.Range("A20").Select()
